$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.95578266666667
$ws.Range("H2").Value = 59.867348
$ws.Range("I2").Value = 0.0117373419656925
$ws.Range("J2").Value = 0.0117373419656925
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 27.81717833333333
$ws.Range("N2").Value = 83.45153500000001
$ws.Range("O2").Value = 0.4044740580248731
$ws.Range("P2").Value = 0.4044740580248732
$ws.Range("Q2").Value = 555.113565219909
$ws.Range("R2").Value = 4996.02208697918
$ws.Range("S2").Value = 0.004747450335289285
$ws.Range("T2").Value = 0.004747450335289286

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.95578266666667
$ws.Range("H3").Value = 59.867348
$ws.Range("I3").Value = 0.0117373419656925
$ws.Range("J3").Value = 0.0117373419656925
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.827410666666666
$ws.Range("N3").Value = 14.482232
$ws.Range("O3").Value = 0.0701926830500802
$ws.Range("P3").Value = 0.0701926830500802
$ws.Range("Q3").Value = 96.33475810674844
$ws.Range("R3").Value = 867.012822960736
$ws.Range("S3").Value = 0.0008238755244482586
$ws.Range("T3").Value = 0.0008238755244482586

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.95578266666667
$ws.Range("H4").Value = 59.867348
$ws.Range("I4").Value = 0.0117373419656925
$ws.Range("J4").Value = 0.0117373419656925
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.06242466666667
$ws.Range("N4").Value = 96.187274
$ws.Range("O4").Value = 0.4662018145637509
$ws.Range("P4").Value = 0.466201814563751
$ws.Range("Q4").Value = 639.8307784143725
$ws.Range("R4").Value = 5758.477005729352
$ws.Range("S4").Value = 0.005471970122561105
$ws.Range("T4").Value = 0.005471970122561106

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.95578266666667
$ws.Range("H5").Value = 59.867348
$ws.Range("I5").Value = 0.0117373419656925
$ws.Range("J5").Value = 0.0117373419656925
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.066688333333334
$ws.Range("N5").Value = 12.200065
$ws.Range("O5").Value = 0.05913144436129575
$ws.Range("P5").Value = 0.05913144436129575
$ws.Range("Q5").Value = 81.15394855306891
$ws.Range("R5").Value = 730.3855369776202
$ws.Range("S5").Value = 0.0006940459833938475
$ws.Range("T5").Value = 0.0006940459833938475

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1637.343343333333
$ws.Range("H6").Value = 4912.03003
$ws.Range("I6").Value = 0.9630320723052701
$ws.Range("J6").Value = 0.9630320723052702
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.81717833333333
$ws.Range("N6").Value = 83.45153500000001
$ws.Range("O6").Value = 0.4044740580248731
$ws.Range("P6").Value = 0.4044740580248732
$ws.Range("Q6").Value = 45546.27177439956
$ws.Range("R6").Value = 409916.4459695961
$ws.Range("S6").Value = 0.3895214902934157
$ws.Range("T6").Value = 0.3895214902934157

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1637.343343333333
$ws.Range("H7").Value = 4912.03003
$ws.Range("I7").Value = 0.9630320723052701
$ws.Range("J7").Value = 0.9630320723052702
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.827410666666666
$ws.Range("N7").Value = 14.482232
$ws.Range("O7").Value = 0.0701926830500802
$ws.Range("P7").Value = 0.0701926830500802
$ws.Range("Q7").Value = 7904.128720602995
$ws.Range("R7").Value = 71137.15848542696
$ws.Range("S7").Value = 0.06759780501838573
$ws.Range("T7").Value = 0.06759780501838575

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1637.343343333333
$ws.Range("H8").Value = 4912.03003
$ws.Range("I8").Value = 0.9630320723052701
$ws.Range("J8").Value = 0.9630320723052702
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 32.06242466666667
$ws.Range("N8").Value = 96.187274
$ws.Range("O8").Value = 0.4662018145637509
$ws.Range("P8").Value = 0.466201814563751
$ws.Range("Q8").Value = 52497.19759909313
$ws.Range("R8").Value = 472474.7783918382
$ws.Range("S8").Value = 0.4489672995918063
$ws.Range("T8").Value = 0.4489672995918064

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1637.343343333333
$ws.Range("H9").Value = 4912.03003
$ws.Range("I9").Value = 0.9630320723052701
$ws.Range("J9").Value = 0.9630320723052702
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.066688333333334
$ws.Range("N9").Value = 12.200065
$ws.Range("O9").Value = 0.05913144436129575
$ws.Range("P9").Value = 0.05913144436129575
$ws.Range("Q9").Value = 6658.565071994663
$ws.Range("R9").Value = 59927.08564795196
$ws.Range("S9").Value = 0.05694547740166243
$ws.Range("T9").Value = 0.05694547740166243

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 17.50081933333334
$ws.Range("H10").Value = 52.502458
$ws.Range("I10").Value = 0.01029341242216722
$ws.Range("J10").Value = 0.01029341242216722
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 27.81717833333333
$ws.Range("N10").Value = 83.45153500000001
$ws.Range("O10").Value = 0.4044740580248731
$ws.Range("P10").Value = 0.4044740580248732
$ws.Range("Q10").Value = 486.8234123747812
$ws.Range("R10").Value = 4381.410711373031
$ws.Range("S10").Value = 0.004163418293317614
$ws.Range("T10").Value = 0.004163418293317614

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 17.50081933333334
$ws.Range("H11").Value = 52.502458
$ws.Range("I11").Value = 0.01029341242216722
$ws.Range("J11").Value = 0.01029341242216722
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.827410666666666
$ws.Range("N11").Value = 14.482232
$ws.Range("O11").Value = 0.0701926830500802
$ws.Range("P11").Value = 0.0701926830500802
$ws.Range("Q11").Value = 84.48364192513957
$ws.Range("R11").Value = 760.352777326256
$ws.Range("S11").Value = 0.0007225222356529418
$ws.Range("T11").Value = 0.0007225222356529419

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 17.50081933333334
$ws.Range("H12").Value = 52.502458
$ws.Range("I12").Value = 0.01029341242216722
$ws.Range("J12").Value = 0.01029341242216722
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 32.06242466666667
$ws.Range("N12").Value = 96.187274
$ws.Range("O12").Value = 0.4662018145637509
$ws.Range("P12").Value = 0.466201814563751
$ws.Range("Q12").Value = 561.1187014799436
$ws.Range("R12").Value = 5050.068313319493
$ws.Range("S12").Value = 0.004798807549267411
$ws.Range("T12").Value = 0.004798807549267413

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 17.50081933333334
$ws.Range("H13").Value = 52.502458
$ws.Range("I13").Value = 0.01029341242216722
$ws.Range("J13").Value = 0.01029341242216722
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.066688333333334
$ws.Range("N13").Value = 12.200065
$ws.Range("O13").Value = 0.05913144436129575
$ws.Range("P13").Value = 0.05913144436129575
$ws.Range("Q13").Value = 71.17037780664114
$ws.Range("R13").Value = 640.5334002597701
$ws.Range("S13").Value = 0.0006086643439292514
$ws.Range("T13").Value = 0.0006086643439292515

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 25.39612333333333
$ws.Range("H14").Value = 76.18836999999999
$ws.Range("I14").Value = 0.01493717330687017
$ws.Range("J14").Value = 0.01493717330687017
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 27.81717833333333
$ws.Range("N14").Value = 83.45153500000001
$ws.Range("O14").Value = 0.4044740580248731
$ws.Range("P14").Value = 0.4044740580248732
$ws.Range("Q14").Value = 706.4484917386611
$ws.Range("R14").Value = 6358.03642564795
$ws.Range("S14").Value = 0.00604169910285059
$ws.Range("T14").Value = 0.006041699102850591

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 25.39612333333333
$ws.Range("H15").Value = 76.18836999999999
$ws.Range("I15").Value = 0.01493717330687017
$ws.Range("J15").Value = 0.01493717330687017
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.827410666666666
$ws.Range("N15").Value = 14.482232
$ws.Range("O15").Value = 0.0701926830500802
$ws.Range("P15").Value = 0.0701926830500802
$ws.Range("Q15").Value = 122.5975166713155
$ws.Range("R15").Value = 1103.37765004184
$ws.Range("S15").Value = 0.001048480271593256
$ws.Range("T15").Value = 0.001048480271593256

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 25.39612333333333
$ws.Range("H16").Value = 76.18836999999999
$ws.Range("I16").Value = 0.01493717330687017
$ws.Range("J16").Value = 0.01493717330687017
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 32.06242466666667
$ws.Range("N16").Value = 96.187274
$ws.Range("O16").Value = 0.4662018145637509
$ws.Range("P16").Value = 0.466201814563751
$ws.Range("Q16").Value = 814.2612912003755
$ws.Range("R16").Value = 7328.35162080338
$ws.Range("S16").Value = 0.006963737300116096
$ws.Range("T16").Value = 0.006963737300116097

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 25.39612333333333
$ws.Range("H17").Value = 76.18836999999999
$ws.Range("I17").Value = 0.01493717330687017
$ws.Range("J17").Value = 0.01493717330687017
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.066688333333334
$ws.Range("N17").Value = 12.200065
$ws.Range("O17").Value = 0.05913144436129575
$ws.Range("P17").Value = 0.05913144436129575
$ws.Range("Q17").Value = 103.2781184715611
$ws.Range("R17").Value = 929.50306624405
$ws.Range("S17").Value = 0.0008832566323102254
$ws.Range("T17").Value = 0.0008832566323102254

